# 20160321_data_sets_V2.xlsx - apply "xlxs with date asstring" edit
#
# Summary of the change (per the authoritative OOXML diff):
#  1. Sheet1 row 9 (the "2016" / Oct-15 DM01 row): the report date in
#     column B moves on from 42248 (2015-09-04) to 42278 (2015-10-04).
#  2. Sheet1 E6 (the stray "wfwqe 10/10/15" note next to the newest row)
#     is removed.
#  3. Sheet1 rows 10-39, column A: the yyyymm numeric period codes
#     (201509, 201508, ... 201304) are replaced with free-text month/year
#     labels ("sep 15", "aug 15", ... "apr 12") - i.e. the period column
#     becomes text instead of a number, matching the commit message
#     "xlxs with date asstring / xlxs con fecha en formato texto".
#  4. The view is scrolled down a bit and the selection moves to A40;
#     columns A, C and D are widened to fit their (new/long) contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- 1. column B: corrected report date on row 9 --------------------------
$ws.Range("B9").Value = 42278

# --- 2. drop the stray note in E6 ------------------------------------------
$ws.Range("E6").ClearContents()

# --- 3. column A, rows 10-39: numeric period -> text month/year label -----
$ws.Range("A10").Value = "sep 15"
$ws.Range("A11").Value = "aug 15"
$ws.Range("A12").Value = "jul 15"
$ws.Range("A13").Value = "jun 15"
$ws.Range("A14").Value = "may 15"
$ws.Range("A15").Value = "apr 15"
$ws.Range("A16").Value = "mar 15"
$ws.Range("A17").Value = "feb 15"
$ws.Range("A18").Value = "jan 15"
$ws.Range("A19").Value = "dec 14"
$ws.Range("A20").Value = "nov 14"
$ws.Range("A21").Value = "oct 14"
$ws.Range("A22").Value = "sep 14"
$ws.Range("A23").Value = "aug 14"
$ws.Range("A24").Value = "jul 14"
$ws.Range("A25").Value = "jun 14"
$ws.Range("A26").Value = "may 14"
$ws.Range("A27").Value = "apr 14"
$ws.Range("A28").Value = "mar 14"
$ws.Range("A29").Value = "feb 14"
$ws.Range("A30").Value = "jan 14"
$ws.Range("A31").Value = "dec 12"
$ws.Range("A32").Value = "nov 12"
$ws.Range("A33").Value = "oct 12"
$ws.Range("A34").Value = "sep 12"
$ws.Range("A35").Value = "aug 12"
$ws.Range("A36").Value = "jul 12"
$ws.Range("A37").Value = "jun 12"
$ws.Range("A38").Value = "may 12"
$ws.Range("A39").Value = "apr 12"

# --- 4. view: scroll so row 5 is at the top, select A40 --------------------
$ws.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A40").Select()

# --- 4b. widen columns A, C, D to fit the new/long content ------------------
$ws.Columns.Item(1).ColumnWidth = 9.7109375
$ws.Columns.Item(3).ColumnWidth = 55.5703125
$ws.Columns.Item(4).ColumnWidth = 146.7109375
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
